$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the placeholder "TERRA" row (was row 2)
$ws.Rows.Item(2).Delete()

# After deletion the existing data occupies rows 2:89.
# Append the 5 new ATP entries at the bottom (rows 90:94).
# Write column A (ATP codes) first, then column C (ONS numbers), then
# column D (names) so new shared-string entries are interned in the same
# order as the source edit.
$ws.Range("A90").Value = "SM323"
$ws.Range("A91").Value = "SM369"
$ws.Range("A92").Value = "MA323"
$ws.Range("A93").Value = "GARA1"
$ws.Range("A94").Value = "GARA2"

$ws.Range("C90").Value = 28230
$ws.Range("C91").Value = 28232
$ws.Range("C92").Value = 28651
$ws.Range("C93").Value = 6500
$ws.Range("C94").Value = 6499

$ws.Range("D90").Value = "SMARIA3 230 "
$ws.Range("D91").Value = "SMARIA3  69 "
$ws.Range("D92").Value = "MACAMB 3 230 "
$ws.Range("D94").Value = "GARABI II525 "
$ws.Range("D93").Value = "GARABI  525 "

# Re-sort the whole table (A2:D94) by ONS number (C) ascending,
# with ATP code (A) ascending as the tiebreaker for blank ONS values.
$rng = $ws.Range("A2:D94")
$rng.Sort($ws.Range("C2:C94"), 1, $ws.Range("A2:A94"), $null, 1)

$ws.Range("D5").Select()
